# Update latest output (run 68)
$wb = $excel.ActiveWorkbook

# --- Sheet: Schedule ---
$ws1 = $wb.Worksheets.Item("Schedule")

$ws1.Range("A2").Value = 46039.27083333334
$ws1.Range("B2").Value = 46039.9375
$ws1.Range("E2").Value = 773.3614980000002
$ws1.Range("F2").Value = 12.78706180555556

$ws1.Range("E3").Value = -69.94760174999999
$ws1.Range("F3").Value = -1.542054712301587

# --- Sheet: Detailed ---
$ws2 = $wb.Worksheets.Item("Detailed")

$ws2.Range("E15").Value = "ON"

$ws2.Range("B45").Value = 9.67895
$ws2.Range("B46").Value = 36.05843
$ws2.Range("B47").Value = 57.06004
$ws2.Range("C47").Value = "historical"
$ws2.Range("E47").Value = "OFF"
$ws2.Range("C48").Value = "historical"
$ws2.Range("C49").Value = "historical"
$ws2.Range("B51").Value = 36.06032
$ws2.Range("B54").Value = 56.97996
$ws2.Range("B55").Value = 56.97996
$ws2.Range("B56").Value = 57.06021
$ws2.Range("B58").Value = 57.06022
$ws2.Range("B59").Value = 57.06003
$ws2.Range("B64").Value = 36.0595
$ws2.Range("B65").Value = 11.42294
$ws2.Range("B67").Value = 0
$ws2.Range("B68").Value = -4.62955
$ws2.Range("B69").Value = -5.50985
$ws2.Range("B70").Value = -1.51945
$ws2.Range("B71").Value = 0
$ws2.Range("B72").Value = -4.89969
$ws2.Range("B73").Value = -0.91435
$ws2.Range("B74").Value = -5.19985
$ws2.Range("B76").Value = -5.51
$ws2.Range("B77").Value = -7
$ws2.Range("B79").Value = -20.85077
$ws2.Range("B80").Value = -14.64697
$ws2.Range("B81").Value = -14.14192
$ws2.Range("B82").Value = -7.19767
$ws2.Range("B83").Value = -6.60177
$ws2.Range("B84").Value = -6.30049
$ws2.Range("B85").Value = 0.2924
$ws2.Range("B86").Value = 2.25724
$ws2.Range("B87").Value = 4.15912
$ws2.Range("B89").Value = 46.36353
$ws2.Range("B90").Value = 57.3
$ws2.Range("B91").Value = 55.89581
$ws2.Range("B92").Value = 56.42265
$ws2.Range("B94").Value = 47.02812
